$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. This shifts the current B:K data
# (and the row-2..20 numeric values) right to C:L, keeping column A in place
# for now (we'll convert it to the numeric index afterwards).
$ws.Range("B:B").Insert()

# Excel's column-insert copies the formatting of the column to the left (A,
# which carries the bold/border "header" style down through row 20) into
# the freshly inserted column. Only the header cell B1 should keep that
# style, so strip the inherited formatting from the data rows B2:B20.
$ws.Range("B2:B20").ClearFormats()

# Give the new header cell B1 the same bold/border/alignment formatting as
# the other header cells (e.g. the now-shifted C1) before setting its text,
# so it reuses the existing header style instead of creating a new one.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New column B header
$ws.Range("B1").Value = "segments"

# Column A currently holds the segment-name strings (with header style s=1
# already applied down A2:A20). Move those names into the new column B,
# then replace column A's values with a 0-based numeric index.
for ($r = 2; $r -le 20; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 1).Value = $r - 2
}
